$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Insert the new "Carrot Cake Recipe" table (with two leading blank
#    paragraphs and one trailing blank paragraph) at the very start of the
#    document body, ahead of the existing "Roasted Butternut Squash" table.
# ---------------------------------------------------------------------------
$introXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblStyle w:val="TableGrid"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="1368"/>
    <w:gridCol w:w="8208"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="1368" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>2019/12/16</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="8208" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Carrot Cake Recipe</w:t>
        </w:r>
        <w:bookmarkStart w:id="0" w:name="_GoBack"/>
        <w:bookmarkEnd w:id="0"/>
      </w:p>
      <w:p>
        <w:r>
          <w:fldChar w:fldCharType="begin"/>
        </w:r>
        <w:r>
          <w:instrText xml:space="preserve"> HYPERLINK "https://www.allrecipes.com/recipe/7402/carrot-cake-iii/" </w:instrText>
        </w:r>
        <w:r>
          <w:fldChar w:fldCharType="separate"/>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rStyle w:val="Hyperlink"/>
          </w:rPr>
          <w:t>https://www.allrecipes.com/recipe/7402/carrot-cake-iii/</w:t>
        </w:r>
        <w:r>
          <w:fldChar w:fldCharType="end"/>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@
$introRange = $d.Range(0, 0)
$introRange.InsertXML($introXml)

# ---------------------------------------------------------------------------
# Helper: walk the document's Paragraphs collection (Next() chaining is the
# reliable way to iterate in this runtime -- Paragraphs.Item(n) and a cell's
# own Range.Paragraphs collection do not index/advance correctly) and find a
# paragraph whose exact text (the visible text, i.e. without the trailing
# paragraph mark, or -- for a paragraph that is also the last one in a table
# cell -- without the trailing paragraph mark + cell mark) matches.  Returns
# $null if not found, otherwise the first match in document order starting
# the search at occurrence index $occurrence (0-based).
# ---------------------------------------------------------------------------
function Find-ParagraphByText($doc, [string]$targetText, [int]$occurrence) {
    $count = $doc.Paragraphs.Count
    $p = $doc.Paragraphs.First
    $seen = 0
    $i = 0
    while ($p -ne $null -and $i -lt $count) {
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $targetText) {
            if ($seen -eq $occurrence) {
                return $p
            }
            $seen = $seen + 1
        }
        $p = $p.Next()
        $i = $i + 1
    }
    return $null
}

# ---------------------------------------------------------------------------
# 2) The first "04/21/2019" table-cell paragraph (Turkey entry) now starts a
#    rendered page, so Word stamps a lastRenderedPageBreak on its run.
# ---------------------------------------------------------------------------
$turkeyDatePara = Find-ParagraphByText $d "04/21/2019" 0
if ($turkeyDatePara -ne $null) {
    $turkeyDateXml = @"
<w:p $wns>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>04/21/2019</w:t>
  </w:r>
</w:p>
"@
    $turkeyDatePara.Range.InsertXML($turkeyDateXml)
} else {
    Write-Output "WARNING: 04/21/2019 paragraph not found"
}

# ---------------------------------------------------------------------------
# 3) In the Cabbage Rolls entry, the page-break mark moves from "Parsley" up
#    to "Tomato puree" (the paragraph immediately above it).
# ---------------------------------------------------------------------------
$tomatoPara = Find-ParagraphByText $d "Tomato puree" 0
if ($tomatoPara -ne $null) {
    $tomatoXml = @"
<w:p $wns>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>Tomato puree</w:t>
  </w:r>
</w:p>
"@
    $tomatoPara.Range.InsertXML($tomatoXml)
} else {
    Write-Output "WARNING: Tomato puree paragraph not found"
}

$parsleyPara = Find-ParagraphByText $d "Parsley" 0
if ($parsleyPara -ne $null) {
    $parsleyXml = @"
<w:p $wns>
  <w:r>
    <w:t>Parsley</w:t>
  </w:r>
</w:p>
"@
    $parsleyPara.Range.InsertXML($parsleyXml)
} else {
    Write-Output "WARNING: Parsley paragraph not found"
}

# ---------------------------------------------------------------------------
# 4) The trailing "_GoBack" bookmark paragraph at the end of the document
#    loses its bookmark (now just an empty paragraph) now that the new table
#    at the top owns the "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastXml = "<w:p $wns/>"
$lastPara.Range.InsertXML($lastXml)

Write-Output "All edits applied"
